$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 157 — this shifts the existing rows 157..193
# down to 158..194 (matching the diff, which shows every row from 157 to 193
# taking on the values that used to belong to the row above it, and a brand
# new row 194 that is a copy of the old row 193).
$ws.Rows.Item(157).Insert()

# Populate the newly inserted row 157 with the new weekly price observation.
$ws.Cells.Item(157, 1).Value = 4
$ws.Cells.Item(157, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(157, 3).Value = "Los Lagos"
$ws.Cells.Item(157, 4).Value = 44551
$ws.Cells.Item(157, 5).Value = 10
$ws.Cells.Item(157, 6).Value = 100112044
$ws.Cells.Item(157, 7).Value = "Perejil"
$ws.Cells.Item(157, 8).Value = "Sin especificar"
$ws.Cells.Item(157, 9).Value = "Primera"
$ws.Cells.Item(157, 10).Value = 160
$ws.Cells.Item(157, 11).Value = 5500
$ws.Cells.Item(157, 12).Value = 5500
$ws.Cells.Item(157, 13).Value = 5500
$ws.Cells.Item(157, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(157, 15).Value = "Región Metropolitana"
$ws.Cells.Item(157, 16).Value = 1833
$ws.Cells.Item(157, 17).Value = 3
$ws.Cells.Item(157, 18).Value = "Hortaliza"
